# Fixed the nan and float issue in friday column
#
# The FRIDAY column (F) held stray ".0"-suffixed numbers (e.g. "2.0",
# "4.0\n...") coming from a pandas float export, plus literal "nan" text
# for the truly-empty bottom rows. Replace them with the clean values
# already used by the rest of the sheet (no trailing ".0") and blank out
# the "nan" cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# F2 held "2.0" and should become plain "2". Column C, row 2 already
# holds that exact clean string ("2") as a shared-string cell, so copy
# it across rather than re-typing a bare numeric literal (which Excel
# would otherwise auto-convert to a *number* instead of keeping it text).
$ws.Range("C2").Copy($ws.Range("F2"))

$ws.Range("F3").Value = "4`nEG-310-04813-Guo-IDE-113A`nEG-209-04812-Beltramo-IDE-107A`nEG-209-04812-Beltramo-IDE-107B"
$ws.Range("F4").Value = "6`nEG-424-06807-Guo-EAX-019`nEG-424-06807-Guo-IDE-107A`nEG-424-06807-Guo-IDE-107B`nEG-360-06811-Sadraey-IDE-204A`nEG-361-06810-Sadique-IDE-118A`nEG-361-06810-Sadique-IDE-118B"
$ws.Range("F5").Value = "8`nEG-419-08810-Rosner-IDE-318`nEG-308-08813-Sadraey-IDE-323`nEG-316-08809-Moghimi-IDE-107A`nEG-316-08809-Moghimi-IDE-107B`nEG-201-08814-TBD-IDE-118A`nEG-201-08814-TBD-IDE-118B`nEG-110-08803-TBD-TBD-TBD`nEG-110-08811-TBD-TBD-TBD"
$ws.Range("F6").Value = "10`nEG-110-10810-Kolenbrander-IDE-128A`nEG-110-10810-Kolenbrander-IDE-128B`nEG-110-10809-Eshed-IDE-128A`nEG-110-10809-Eshed-IDE-128B`nEG-335-10811-Daigneau-IDE-206A`nEG-335-10811-Daigneau-IDE-206B`nEG-409-10812-TBD-TBD-TBD"
$ws.Range("F7").Value = "12`nEG-410-12808-Sadique-IDE-318`nEG-316-12806-TBD-IDE-118A`nEG-316-12806-TBD-IDE-118B"

# F8:F10 only ever held the literal text "nan" -- clear them out entirely.
$ws.Range("F8:F10").ClearContents()
